$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.515.02'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.516.97'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '541.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.77'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.83%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.563'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.89%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.524.03'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.83%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.101'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.49%  '
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.43'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.36%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.352'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.966.55'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.33'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '59.427.21'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000141'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.53%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.512.18'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.07'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.69'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.37%  '
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.88'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.10'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.421'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.168'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('E27').Value = '  +0.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.78'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0784'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.75'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.79'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '165.56'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.11'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -10.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.39'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.98%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.54'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.19'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.59'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.67'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.93%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.814'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.27'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -9.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '281.32'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.997'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.88'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.598'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '126.23'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0939'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0513'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0223'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.85'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.77%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.774.69'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.70%  '
